# Daily attendance processing - swap order of "Recorded By" names
# Change "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# in every row of column G ("Recorded By") on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Determine the last used row in column A (data starts at row 2, row 1 is the header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
